$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("configuration")

# Change the active stage from "qa" to "prod" so the VLOOKUP picks up the
# production frontend URL for the new login test case.
$ws.Range("A2").Value = "prod"

$ws.Activate()
$ws.Range("C4").Select()
